$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters used in this data block (H is intentionally skipped - stays empty)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New vm_pu values per row (row index => values for columns B,C,D,E,F,I,J,K,L,M,N)
$data = @{
    2 = @(1.02,1.022213354495531,1.031946749316871,0.9926147277508489,1.038518555771472,1.030219767504432,1.027400456009924,1.034753398002257,0.9955398523336033,1.04130635982613,1.013195684962477)
    3 = @(1.02,1.023049794871936,1.032595775572389,0.9936372048519304,1.039351025309873,1.030324312505385,1.027875436477408,1.035211521753906,0.9963617723202692,1.041948799947362,1.013353864146786)
    4 = @(1.02,1.023591573128362,1.033016024765011,0.9942998659930995,1.039890487155654,1.030390683779427,1.028182686297131,1.035507559126296,0.9968940712668345,1.042364652861201,1.013456150574995)
    5 = @(1.02,1.023819465725599,1.033192764124587,0.9945786998346017,1.040117466112167,1.030418280037187,1.02831183031977,1.035631916262068,0.997117960005301,1.042539511973366,1.013499135525952)
    6 = @(1.02,1.023857737416069,1.033222443293066,0.9946255319796338,1.040155587915889,1.030422895593609,1.028333512738065,1.035652790645647,0.9971555583673453,1.042568873565596,1.013506351921758)
    7 = @(1.02,1.023594617735389,1.033018386105593,0.9943035907982488,1.039893519318606,1.03039105372607,1.028184412020945,1.035509221175209,0.9968970624462087,1.042366989203544,1.013456725006476)
    8 = @(1.02,1.022495919574101,1.032166030505041,0.9929600610674301,1.038799726272138,1.030255362847524,1.027560996672622,1.034908304926108,0.995817528259106,1.041523443175166,1.013249155774469)
    9 = @(1.02,1.020564118804924,1.030666340480161,0.9906006454969559,1.036878523570111,1.03000651258276,1.026461791031764,1.033846416075593,0.9939188001724441,1.040038238151793,1.012882905741685)
    10 = @(1.02,1.019279198386136,1.029668181514221,0.989033133672735,1.035602006003683,1.029834099849149,1.025728610852808,1.033136556525253,0.9926553831429383,1.039049027166638,1.012638438288295)
    11 = @(1.02,1.018723531894759,1.029236378211678,0.988355674866747,1.035050300293431,1.029757906909498,1.025411061766187,1.032828738084476,0.9921088820399291,1.038620928117932,1.012532515386208)
    12 = @(1.02,1.018517241207839,1.029076050101707,0.9881042295826724,1.034845529661271,1.029729375184482,1.025293099518384,1.032714335196063,0.9919059725120875,1.038461950162705,1.01249316129972)
    13 = @(1.02,1.018561486321827,1.029110438151237,0.9881581567098651,1.034889446501282,1.029735505746243,1.025318403270945,1.032738877941986,0.9919494934313052,1.038496049768594,1.012501603312365)
    14 = @(1.02,1.018706477618129,1.029223124147039,0.9883348863814464,1.035033370666157,1.02975555316195,1.025401311174615,1.032819282832492,0.9920921077337197,1.038607786184282,1.012529262558506)
    15 = @(1.02,1.018795825976293,1.02929256209329,0.9884438009545853,1.035122068014377,1.029767874542389,1.02545239213179,1.032868814312441,0.9921799884222134,1.038676635642536,1.012546303073506)
    16 = @(1.02,1.019316091284655,1.029696847591553,0.9890781214508737,1.035638642864761,1.029839124194871,1.025749684017273,1.033156976172441,0.9926916645766087,1.039077443801183,1.012645466671234)
    17 = @(1.02,1.01964263181773,1.02995055523829,0.989476357848556,1.03596295505779,1.029883406159872,1.025936147550756,1.03333761444297,0.9930127773699352,1.039328924482752,1.012707651783602)
    18 = @(1.02,1.019833166031397,1.030098577676838,0.9897087662937556,1.036152220586567,1.029909086730166,1.026044901041361,1.03344293470302,0.9932001317071769,1.039475631553002,1.012743916824736)
    19 = @(1.02,1.019898144871051,1.030149056072217,0.9897880325774034,1.036216772094209,1.029917817958206,1.026081981856082,1.033478838856697,0.9932640239640975,1.039525658662872,1.012756281153481)
    20 = @(1.02,1.019607589997447,1.029923330755294,0.9894336180360679,1.035928149100726,1.029878670463985,1.025916142548692,1.03331823810248,0.9929783193494215,1.039301940631501,1.012700980576903)
    21 = @(1.02,1.018663778285216,1.029189939179759,0.9882828385668249,1.034990984225127,1.029749656049138,1.025376897140065,1.032795607390963,0.9920501090198102,1.038574881557551,1.012521117861475)
    22 = @(1.02,1.018070995428562,1.028729191404153,0.9875604150241495,1.034402663303341,1.029667207592021,1.025037792710032,1.032466630935827,0.9914670000341481,1.038117966088229,1.012407975684464)
    23 = @(1.02,1.018385180697488,1.028973407340015,0.9879432794643023,1.034714456269089,1.029711041115917,1.025217563609786,1.03264106292897,0.991776070289318,1.038360164649375,1.012467959591659)
    24 = @(1.02,1.019623423675746,1.029935632206835,0.9894529299347244,1.035943876106661,1.029880810779566,1.025925181971443,1.033326993569689,0.9929938892766442,1.039314133402611,1.012703995028136)
    25 = @(1.02,1.021063022976906,1.031053766161843,0.9912096547607049,1.037374453840564,1.030071997141501,1.026746033711423,1.034121286598878,0.9944092447426414,1.040422043270347,1.012977645092706)
}

foreach ($r in $data.Keys) {
    $rowValues = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $ws.Range("$colLetter$r").Value = $rowValues[$i]
    }
}
